$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventario")

# --- Update existing rows (Stock/column D adjustments, etc.) ---

# Row 2: Stock 23 -> 21
$ws.Range("D2").Value = 21.0

# Row 3: Stock 7 -> 0
$ws.Range("D3").Value = 0.0

# Row 5: Stock 28 -> 23
$ws.Range("D5").Value = 23.0

# Row 7: product renamed "Talco para bebé" -> "Talco para pies",
# price 30 -> 25, stock 8 -> 5
$ws.Range("B7").Value = "Talco para pies"
$ws.Range("C7").Value = 25.0
$ws.Range("D7").Value = 5.0

# Row 9: Stock 17 -> 11
$ws.Range("D9").Value = 11.0

# --- Remove the "2452 / Bio-Electro 24 tabletas" and ---
# --- "9498 / Almetec 40mg 28 tabletas" rows (old rows 11-12); ---
# --- this shifts old rows 13-17 up to become rows 11-15 ---
$ws.Range("A11:D12").EntireRow.Delete()

# After the shift, update the Stock values that differ from a
# straight carry-up of the old data
$ws.Range("D11").Value = 17.0
$ws.Range("D13").Value = 43.0
$ws.Range("D14").Value = 5.0

# --- Append the new product row (row 16) ---
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "6771"
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").Value = "Aspirina 250mg"
$ws.Range("C16").Value = 35.0
$ws.Range("D16").Value = 20.0
